$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values (these feed the formulas in D2/E2/E3/F3 which recalc automatically)
$ws.Range("E6").Value = 0.95
$ws.Range("E8").Value = 0.95
$ws.Range("E9").Value = 0.9

# Scroll the sheet view so that row 3 becomes the top-left visible cell
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
